$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the entire first row (as a user would before deleting it), then delete it.
# This mirrors deleting a header row: all rows below shift up by one.
$ws.Rows.Item(1).Select()
$ws.Rows.Item(1).Delete()
